$d = $word.ActiveDocument

# 1. Remove the paragraph "EasyStar(TM) + ablador com forca de contato e sistema EverPace."
#    (the whole paragraph, including its paragraph mark, is deleted).
$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*ablador com for*a de contato e sistema EverPace*") {
        $p.Range.Delete()
        $found = $true
    }
}

# 2. Prefix every "Materiais:" list item with a bullet glyph + space.
$items = @(
    "Cateter de Navegação – EasyStar™ 3D Sensor",
    "Conectores EasyStar™ (cabos de localização e sistema)",
    "Cateter de Ablação com Força de Contato – FireMagic™ 3D Plus + conector",
    "Cateter Decapolar – EasyFinder™ Deca + conector",
    "Patch de Impedância – EverPace™ Impedance Patch",
    "Patch de Referência – EverPace™ Reference Patch",
    "Bainha Curva Fixa – EasySheath™ Fixed Curve",
    "Bainha Defletiva – SureFlex™ Steering Introducer",
    "Agulha de Punção",
    "Introdutor – 3"
)

foreach ($item in $items) {
    $rng = $d.Content
    $matched = $rng.Find.Execute($item, $true, $false, $false, $false, $false, $true, 1, $false, "• " + $item, 2)
}
